$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").Value = "2024-10-11 00:00:00"
$ws.Range("B71").Value = 76450
$ws.Range("C71").Value = 10776.11
$ws.Range("D71").Value = 9536.379999999999
$ws.Range("E71").Value = 7.0673
